$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 998
$ws.Range("J17").Value = 998
$ws.Range("L17").Value = 2994
$ws.Range("N17").Value = -3330

$ws.Range("H64").Value = 4632.6665
$ws.Range("J64").Value = 4449.5
$ws.Range("L64").Value = 4449.5
$ws.Range("N64").Value = -4945.5

$ws.Range("H67").Value = 4632.6665
$ws.Range("J67").Value = 4449.5
$ws.Range("L67").Value = 4449.5
$ws.Range("N67").Value = -6165.5

$ws.Range("H70").Value = 1447.1666
$ws.Range("I70").Value = 1371.3572
$ws.Range("J70").Value = 1712.5
$ws.Range("K70").Value = 4114.071599999999
$ws.Range("L70").Value = 5137.5
$ws.Range("M70").Value = -3844.071599999999
$ws.Range("N70").Value = -5677.5

$ws.Range("H73").Value = 1447.1666
$ws.Range("I73").Value = 1371.3572
$ws.Range("J73").Value = 1712.5
$ws.Range("K73").Value = 4114.071599999999
$ws.Range("L73").Value = 5137.5
$ws.Range("M73").Value = -3178.071599999999
$ws.Range("N73").Value = -7009.5

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H76").Value = 9000
$ws.Range("I76").Value = 8000
$ws.Range("K76").Value = 8000
$ws.Range("M76").Value = -7685

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H79").Value = 9000
$ws.Range("I79").Value = 8000
$ws.Range("K79").Value = 8000
$ws.Range("M79").Value = -6908

$ws.Range("H92").Value = 2040
$ws.Range("I92").Value = 1925
$ws.Range("J92").Value = 2500
$ws.Range("K92").Value = 1925
$ws.Range("L92").Value = 2500
$ws.Range("M92").Value = -677
$ws.Range("N92").Value = -4996

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()

$ws.Range("H116").Value = 10555.2
$ws.Range("I116").Value = 5694
$ws.Range("K116").Value = 5694
$ws.Range("M116").Value = -2252

$ws.Range("H125").Value = 2900
$ws.Range("J125").Value = 2900
$ws.Range("L125").Value = 26100
$ws.Range("N125").Value = -31020

$ws.Range("H138").Value = 5796.3335
$ws.Range("J138").Value = 5796.3335
$ws.Range("L138").Value = 17389.0005
$ws.Range("N138").Value = -27669.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 226.875
$ws.Range("I97").Value = 146
$ws.Range("J97").Value = 361.66666
$ws.Range("K97").Value = 146
$ws.Range("L97").Value = 361.66666
$ws.Range("M97").Value = 350
$ws.Range("N97").Value = -1353.66666

$ws.Range("H132").Value = 9152
$ws.Range("I132").Value = 9152
$ws.Range("K132").Value = 27456
$ws.Range("M132").Value = -24926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1712.5
$ws.Range("I107").Value = 2075
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 2075
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = -155
$ws.Range("N107").Value = -5190

$ws.Range("H134").Value = 3078
$ws.Range("I134").Value = 1996.6666
$ws.Range("J134").Value = 4700
$ws.Range("K134").Value = 5989.9998
$ws.Range("L134").Value = 14100
$ws.Range("M134").Value = -3454.9998
$ws.Range("N134").Value = -19170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 41349.75
$ws.Range("J43").Value = 41349.75
$ws.Range("L43").Value = 41349.75
$ws.Range("N43").Value = -41717.75

$ws.Range("H96").Value = 18849.7
$ws.Range("J96").Value = 18849.7
$ws.Range("L96").Value = 18849.7
$ws.Range("N96").Value = -24341.7

$ws.Range("H101").Value = 41349.75
$ws.Range("J101").Value = 41349.75
$ws.Range("L101").Value = 41349.75
$ws.Range("N101").Value = -47839.75

$ws.Range("H102").Value = 38498.332
$ws.Range("J102").Value = 38498.332
$ws.Range("L102").Value = 38498.332
$ws.Range("N102").Value = -43366.332

$ws.Range("H103").Value = 37499.5
$ws.Range("I103").Value = 37499.5
$ws.Range("K103").Value = 37499.5
$ws.Range("M103").Value = -36327.5

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 750
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 750
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 2250
$ws.Range("N26").Value = -2826
$ws.Range("M26").ClearContents()

$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6372
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 13500
$ws.Range("N83").Value = -22860
$ws.Range("M83").ClearContents()

$ws.Range("H117").Value = 3123.5
$ws.Range("J117").Value = 5374.6665
$ws.Range("L117").Value = 16123.9995
$ws.Range("N117").Value = -23007.9995

$ws.Range("H122").Value = 697.5
$ws.Range("I122").Value = 697.5
$ws.Range("K122").Value = 6277.5
$ws.Range("M122").Value = -3827.5

$ws.Range("H129").Value = 1643.6666
$ws.Range("I129").Value = 1025.8
$ws.Range("J129").Value = 4733
$ws.Range("K129").Value = 3077.4
$ws.Range("L129").Value = 14199
$ws.Range("M129").Value = 1922.6
$ws.Range("N129").Value = -24199

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2825
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2825
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2825
$ws.Range("N7").Value = -3049
$ws.Range("M7").ClearContents()

$ws.Range("H18").Value = 1868
$ws.Range("I18").Value = 1868
$ws.Range("K18").Value = 1868
$ws.Range("M18").Value = -1696

$ws.Range("H20").Value = 44999.75
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1996.6666
$ws.Range("I82").Value = 1995
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1995
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -1634
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 1996.6666
$ws.Range("I85").Value = 1995
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1995
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -747
$ws.Range("N85").Value = -4496

$ws.Range("H94").Value = 64748.5
$ws.Range("J94").Value = 64748.5
$ws.Range("L94").Value = 64748.5
$ws.Range("N94").Value = -66100.5

$ws.Range("H126").Value = 2825
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2825
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8475
$ws.Range("N126").Value = -13415
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

$ws.Range("H122").Value = 1134.25
$ws.Range("I122").Value = 1151.2
$ws.Range("J122").Value = 1049.5
$ws.Range("K122").Value = 3453.6
$ws.Range("L122").Value = 3148.5
$ws.Range("M122").Value = -1003.6
$ws.Range("N122").Value = -8048.5

$ws.Range("H126").Value = 2998.25
$ws.Range("I126").Value = 2998.25
$ws.Range("K126").Value = 8994.75
$ws.Range("M126").Value = -6524.75

$ws.Range("H132").Value = 1656.5385
$ws.Range("I132").Value = 1516.1111
$ws.Range("K132").Value = 4548.3333
$ws.Range("M132").Value = -2018.3333

$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280
